$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (changed/updated) date column C for rows 2-7
# from serial 45170 (2023-09-01) to serial 45174 (2023-09-05)
foreach ($r in 2..7) {
    $ws.Cells.Item($r, 3).Value = 45174
}
